# Atualização de bases das ligas, do dia: 25-04-2024 às 21:26
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 73 and 74: swap all match data (columns B, E:AB) between the two rows ---
# Row 73 (before) values, columns B, E:AB
$row73 = @{
    B = 7646750
    E = "Perth Glory"
    F = "Wellington Phoenix"
    G = 3
    H = 4
    I = "A"
    J = 2.45
    K = 3.75
    L = 2.55
    M = 3.1
    N = 3.8
    O = 2.05
    P = 0.25
    Q = 2
    R = 1.85
    S = 3
    T = 1.925
    U = 1.925
    V = -1
    W = -1
    X = 1.05
    Y = -1
    Z = 0.8500000000000001
    AA = 0.925
    AB = -1
}

# Row 74 (before) values, columns B, E:AB
$row74 = @{
    B = 7646749
    E = "Brisbane Roar"
    F = "Newcastle Jets"
    G = 3
    H = 2
    I = "H"
    J = 1.909
    K = 4
    L = 3.4
    M = 2.4
    N = 4
    O = 2.6
    P = 0
    Q = 1.83
    R = 2.07
    S = 3.25
    T = 1.9
    U = 1.95
    V = 1.4
    W = -1
    X = -1
    Y = 0.8300000000000001
    Z = -1
    AA = 0.8999999999999999
    AB = -1
}

foreach ($col in $row74.Keys) {
    $ws.Range(($col + "73")).Value = $row74[$col]
}
foreach ($col in $row73.Keys) {
    $ws.Range(($col + "74")).Value = $row73[$col]
}

# Row 112 (before) values, columns B, E:AB
$row112 = @{
    B = 7127376
    E = "Newcastle Jets"
    F = "Macarthur FC"
    G = 2
    H = 2
    I = "D"
    J = 1.95
    K = 4
    L = 3.4
    M = 1.909
    N = 4.2
    O = 3.6
    P = -0.5
    Q = 1.89
    R = 2.01
    S = 3.5
    T = 1.95
    U = 1.9
    V = -1
    W = 3.2
    X = -1
    Y = -1
    Z = 1.01
    AA = 0.95
    AB = -1
}

# Row 113 (before) values, columns B, E:AB
$row113 = @{
    B = 7127379
    E = "Melbourne Victory"
    F = "Central Coast Mariners"
    G = 0
    H = 1
    I = "A"
    J = 1.95
    K = 3.6
    L = 3.8
    M = 1.909
    N = 3.6
    O = 4
    P = -0.5
    Q = 1.9
    R = 1.95
    S = 2.75
    T = 1.925
    U = 1.925
    V = -1
    W = -1
    X = 3
    Y = -1
    Z = 0.95
    AA = -1
    AB = 0.925
}

foreach ($col in $row113.Keys) {
    $ws.Range(($col + "112")).Value = $row113[$col]
}
foreach ($col in $row112.Keys) {
    $ws.Range(($col + "113")).Value = $row112[$col]
}

# --- Rows 158-163: update specific odds values ---
$ws.Range("M158").Value = 2.1
$ws.Range("O158").Value = 3
$ws.Range("Q158").Value = 1.9
$ws.Range("R158").Value = 2
$ws.Range("S158").Value = 3.75
$ws.Range("T158").Value = 1.95
$ws.Range("U158").Value = 1.9

$ws.Range("M159").Value = 3.8
$ws.Range("N159").Value = 4
$ws.Range("O159").Value = 1.85
$ws.Range("Q159").Value = 2.02
$ws.Range("R159").Value = 1.88
$ws.Range("T159").Value = 1.95
$ws.Range("U159").Value = 1.9

$ws.Range("M160").Value = 1.666
$ws.Range("N160").Value = 4.333
$ws.Range("O160").Value = 4.333
$ws.Range("Q160").Value = 1.88
$ws.Range("R160").Value = 2.02
$ws.Range("S160").Value = 3.25
$ws.Range("T160").Value = 2.05
$ws.Range("U160").Value = 1.8

$ws.Range("M161").Value = 1.8
$ws.Range("N161").Value = 3.9
$ws.Range("O161").Value = 4.2
$ws.Range("Q161").Value = 1.87
$ws.Range("R161").Value = 2.03
$ws.Range("T161").Value = 1.9
$ws.Range("U161").Value = 1.95

$ws.Range("N162").Value = 5
$ws.Range("O162").Value = 7
$ws.Range("Q162").Value = 1.99
$ws.Range("R162").Value = 1.91
$ws.Range("S162").Value = 3.75
$ws.Range("T162").Value = 2
$ws.Range("U162").Value = 1.85

$ws.Range("N163").Value = 5.25
$ws.Range("O163").Value = 6
$ws.Range("T163").Value = 1.975
$ws.Range("U163").Value = 1.875

